$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.929.68'
$ws.Range('E2').Value = '  +2.04%  '
$ws.Range('D3').Value = '2.330.99'
$ws.Range('E3').Value = '  -0.56%  '
$ws.Range('D4').Value = "'0.998"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = "'544.17"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.93%  '
$ws.Range('D6').Value = "'134.38"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('D7').Value = "'0.995"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').Value = "'0.537"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').Value = '2.360.66'
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('D11').Value = "'0.154"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').Value = "'5.39"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = "'0.356"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.39%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = "'23.62"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.70%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.750.48'
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('D16').Value = '57.684.52'
$ws.Range('E16').Value = '  +1.59%  '
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').Value = '2.346.52'
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('D19').Value = "'337.12"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.39%  '
$ws.Range('D20').Value = "'10.50"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').Value = "'4.22"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').Value = "'6.74"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('D23').Value = "'0.995"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').Value = "'62.08"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.07%  '
$ws.Range('D25').Value = "'0.169"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.48%  '
$ws.Range('D26').Value = "'8.54"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.18%  '
$ws.Range('D27').Value = "'0.992"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.58%  '
$ws.Range('D28').Value = "'1.37"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.24%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = "'1.77"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.33%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = "'173.61"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.63%  '
$ws.Range('D31').Value = '0.0₃0741'
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('D32').Value = "'6.17"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.23%  '
$ws.Range('D33').Value = "'18.52"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('B34').Value = 'SuiNetwork'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D34').Value = "'1.01"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +13.68%  '
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').Value = "'0.998"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('D38').Value = "'4.08"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('D40').Value = "'39.28"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('D41').Value = "'149.71"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('D42').Value = "'0.377"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.12%  '
$ws.Range('D43').Value = "'3.64"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.73%  '
$ws.Range('D44').Value = "'285.68"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.51%  '
$ws.Range('D45').Value = "'0.0932"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = "'0.0505"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.48%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = "'0.562"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = "'18.78"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.77%  '
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'17.58"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.56%  '
$ws.Range('B51').Value = 'Polygon'
$ws.Range('C51').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D51').Value = "'0.382"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +9.73%  '
